# Insert a new data row into the "Vega Modelo de Temuco - Betarraga" sheet.
# The new record is inserted at row 320, pushing all subsequent rows (old
# 320..399) down by one (to 321..400) and extending the used range from
# A1:R399 to A1:R400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 320..399 down by one row, creating a blank row 320.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new record's values.
$ws.Cells.Item(320, 1).Value = 10
$ws.Cells.Item(320, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(320, 3).Value = "La Araucanía"
$ws.Cells.Item(320, 4).Value = 44736
$ws.Cells.Item(320, 4).NumberFormat = $ws.Cells.Item(321, 4).NumberFormat
$ws.Cells.Item(320, 5).Value = 9
$ws.Cells.Item(320, 6).Value = 100114014
$ws.Cells.Item(320, 7).Value = "Betarraga"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 50
$ws.Cells.Item(320, 11).Value = 8000
$ws.Cells.Item(320, 12).Value = 8000
$ws.Cells.Item(320, 13).Value = 8000
$ws.Cells.Item(320, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(320, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(320, 16).Value = 667
$ws.Cells.Item(320, 17).Value = 12
$ws.Cells.Item(320, 18).Value = "Hortaliza"
